# Rename the template worksheet from the old "export to production" label to
# the new "import order" (don nhap) label, and move the active selection to
# B32 to match the author's final cursor position when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Template đơn nhập"

$ws.Range("B32").Select()
